$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting B:E (previously A:D) to the right.
$ws.Columns.Item(1).Insert()

# Fix up header row text: "(n=5080)" -> "(n = 5080)" for the two headers that moved to D1/E1.
$ws.Range("D1").Value = "Treatment at T1 (n = 5080)"
$ws.Range("E1").Value = "Treatment at T2 (n = 5080)"

# Fix up the category labels that moved from A3/A10 to B3/B10, adding spaces around "=".
$ws.Range("B3").Value = "Gender (P = 1.000)"
$ws.Range("B10").Value = "Interested in News (P = 1.000)"

# Apply the header style to the new column A cells A2:A15 (A1 stays blank/unstyled).
$ws.Range("D1").Copy()
$ws.Range("A2:A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
